$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-code" and "codeforiati:group-name" columns (C and D)
# need to be swapped, including the header row, across all used rows.
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value()
    $dVal = $dCell.Value()
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
